$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "25 uur"

$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "32 uur"

$ws.Range("A9").Value = 13
$ws.Range("B9").Value = "45 uur"

$ws.Range("B10").Select()
